# Update the Status column of the Product backlog table.
# Rows (by ID) 01 and 02: "In progress" -> "Done"
# Rows (by ID) 03, 04, 05 and 06: "To be started" -> "In progress"
# Rows 07, 08, 09 remain "To be started" (untouched)

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Status column is the last (7th) column; data rows are 2..10 (row 1 is the header)
$statusCol = 7

$newStatus = @{
    2  = "Done"          # ID 01
    3  = "Done"          # ID 02
    4  = "In progress"   # ID 03
    5  = "In progress"   # ID 04
    6  = "In progress"   # ID 05
    7  = "In progress"   # ID 06
}

foreach ($rowIndex in $newStatus.Keys) {
    $cell = $table.Cell($rowIndex, $statusCol)
    $cellRange = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters before replacing text
    $cellRange.End = $cellRange.End - 1
    $cellRange.Text = $newStatus[$rowIndex]
}
